# Reorders the "Recorded By" list in column G so the most recently
# appended recorder (the last comma-separated entry) is shown first,
# e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com".
# Rows that already start with "System" are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val.Split(",")

        if ($parts.Length -gt 1) {
            $trimmed = @()
            foreach ($p in $parts) {
                $trimmed += $p.Trim()
            }

            if (-not $trimmed[0].Equals("System")) {
                $lastIdx = $trimmed.Length - 1
                $lastItem = $trimmed[$lastIdx]
                $rest = $trimmed[0..($lastIdx - 1)]
                $newParts = @($lastItem) + $rest
                $result = [string]::Join(", ", $newParts)
                $cell.Value = $result
            }
        }
    }
}
